$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D), Volume(1h) (E) and Hora (G) columns store their values as
# literal text (e.g. "301.24", "-0.64%", "5"), matching the rest of the sheet.
# Force a Text number format on each edited cell before assigning its value so
# Excel does not auto-convert these numeric-looking strings into real numbers
# or percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.64%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.97%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "5"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.978"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.52%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "5"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07735"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.38%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.205"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.91%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "5"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.006"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.08%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "5"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.988"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.05%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "5"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.10%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "5"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09344"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.31%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "5"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1797"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.25%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "5"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08411"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.87%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "5"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03528"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.25%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "5"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09917"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "5"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001467"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.28%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "5"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005698"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.77%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "5"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.475"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.39%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "5"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.054"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.14%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "5"

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "5"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1314"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.50%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "5"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.573"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.20%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "5"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2227"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.19%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "5"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04648"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.77%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "5"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001227"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.12%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "5"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004437"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.74%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "5"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.20%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "5"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004742"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "39.76%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "5"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "5"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "5"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "5"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "5"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "5"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "5"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "5"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "5"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "5"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "5"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "5"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01752"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.02%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "5"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04689"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.43%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "5"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007931"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.56%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "5"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.75%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "5"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007659"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.30%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "5"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002286"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.80%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "5"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01019"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.11%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "5"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006107"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.10%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "5"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.18%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "5"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.662"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "182.82%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "5"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "5"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.18%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "5"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.18%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "5"

